# Append a "Late Entry" follow-up note at the end of the log, after the
# existing "(Note to self: ...)" paragraph, separated by one blank
# paragraph. The trailing _GoBack bookmark (which marks the last edit
# position) should end up wrapping the newly typed text, exactly as it
# would if a person had put their cursor at the end of the document and
# typed the new content.

$d = $word.ActiveDocument
$cr = [char]13

# Sanity check: confirm the document really ends with the expected note.
$lastPara = $d.Paragraphs.Last
if ($lastPara.Range.Text -notlike "*Check in the morning, and report on findings)*") {
    throw "Unexpected document tail, aborting: [$($lastPara.Range.Text)]"
}

$newText = "Late Entry: Found the problem in the Matrix and corrected it. Just had things linked up incorrectly."

# Insert a unique two-character marker plus the new sentence at the very
# end of the document, then surgically turn the two marker characters
# into paragraph marks. Doing the split this way (rather than via
# InsertParagraphAfter/Before, or by embedding the carriage returns
# directly in the InsertBefore call) avoids leaving a stray empty run
# behind and keeps the relocated bookmark ordered after the new run,
# matching how Word itself records this kind of edit.
$marker = "@@"
$insertPos = $d.Content.End
$rng = $d.Range($insertPos, $insertPos)
$rng.InsertBefore($marker + $newText)

$m1 = $d.Range($insertPos - 1, $insertPos)
if ($m1.Text -ne "@") { throw "marker 1 mismatch: [$($m1.Text)]" }
$m1.Text = "" + $cr

$m2 = $d.Range($insertPos, $insertPos + 1)
if ($m2.Text -ne "@") { throw "marker 2 mismatch: [$($m2.Text)]" }
$m2.Text = "" + $cr

Write-Output "Inserted late entry; paragraph count now $($d.Paragraphs.Count)"
